$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (before the current K_SEX row), shifting
# K_SEX and K_URBAN down by one row.
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the K_SERIES entry.
$ws.Range("A9").Value = "K_SERIES"
$ws.Range("B9").Value = "Zeitreihe"
$ws.Range("C9").Value = "Time series"

# Carry over the data-row formatting (font/fill/border/alignment) from the
# row below, matching the look of every other category row.
$ws.Range("A10:C10").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
